# This workbook contains 28 worksheets (one per backward-elimination step),
# each holding a statsmodels OLS summary text dump in cell B2. The summary
# text embeds the date/time the regression was run:
#   Date:                Sat, 28 Dec 2019   ...
#   Time:                        20:59:52   ...
# The source run was regenerated a day later, so every sheet's Date: line
# moves from "Sat, 28 Dec 2019" to "Sun, 29 Dec 2019", and the Time: line
# is updated to the new run time (16:11:22 for the first 15 sheets,
# 16:11:23 for the remaining 13 sheets, reflecting the wall-clock tick
# that occurred partway through the batch run).

$wb = $excel.ActiveWorkbook

$oldDate = "Sat, 28 Dec 2019"
$newDate = "Sun, 29 Dec 2019"
$oldTime = "20:59:52"

$sheetCount = $wb.Worksheets.Count

for ($i = 1; $i -le $sheetCount; $i++) {
    $ws = $wb.Worksheets.Item($i)
    $cell = $ws.Range("B2")
    $text = $cell.Value2

    if ($i -le 15) {
        $newTime = "16:11:22"
    } else {
        $newTime = "16:11:23"
    }

    $text = $text.Replace($oldDate, $newDate)
    $text = $text.Replace($oldTime, $newTime)

    $cell.Value2 = $text
}
